# Split the first paragraph's sentence into 4 runs by appending
# " (Changed main)" as three additional runs after the existing
# "This is a Microsoft word document." run, matching:
#   <w:r><w:t>This is a Microsoft word document.</w:t></w:r>
#   <w:r><w:t xml:space="preserve"> (</w:t></w:r>
#   <w:r><w:t>Changed main</w:t></w:r>
#   <w:r><w:t>)</w:t></w:r>
#
# Word COM doesn't expose a "Run" object directly, so the reliable way
# to force distinct <w:r> elements (even though they share identical
# run formatting, which would otherwise get coalesced into a single
# run on save) is to replace the paragraph's text range with an
# equivalent WordprocessingML fragment via Range.InsertXML - the same
# mechanism Word itself uses under the hood for OOXML-fragment pastes.

$d = $word.ActiveDocument
$p1 = $d.Paragraphs(1)

$r = $p1.Range
# Exclude the trailing paragraph mark from the range so InsertXML
# replaces only the run content, leaving the paragraph mark (and its
# paragraph-level properties) untouched.
$r.End = $r.End - 1

$wordNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$xml = '<?xml version="1.0" standalone="yes"?>' +
       '<?mso-application progid="Word.Document"?>' +
       '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
         '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData>' +
             '<w:document ' + $wordNs + '>' +
               '<w:body>' +
                 '<w:p>' +
                   '<w:r><w:t>This is a Microsoft word document.</w:t></w:r>' +
                   '<w:r><w:t xml:space="preserve"> (</w:t></w:r>' +
                   '<w:r><w:t>Changed main</w:t></w:r>' +
                   '<w:r><w:t>)</w:t></w:r>' +
                 '</w:p>' +
               '</w:body>' +
             '</w:document>' +
           '</pkg:xmlData>' +
         '</pkg:part>' +
       '</pkg:package>'

$r.InsertXML($xml)
